$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9579288025889967
$ws.Range("C2").Value = 0.8058252427184466

$ws.Range("B3").Value = 0.9676375404530745
$ws.Range("C3").Value = 0.8381877022653722

$ws.Range("B4").Value = 0.9741100323624595
$ws.Range("C4").Value = 0.7702265372168284

$ws.Range("B5").Value = 0.9676375404530745
$ws.Range("C5").Value = 0.8284789644012945

$ws.Range("B6").Value = 0.9644012944983819
$ws.Range("C6").Value = 0.8058252427184466
